$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.45114143204246204,
    -0.0099999989407990597,
    -0.0089999989204816444,
    -0.011999999737678735,
    -0.0059999989307222279,
    -0.0059999988961259021,
    -0.019999998720535928,
    -0.019999998714102851,
    -0.0059999988816459293,
    -0.0059999988776908708,
    -0.0044999988962821647,
    -0.0059999988771664015,
    -0.0059999988765184753,
    -0.011999998802836309,
    -0.0059999988787842184,
    -0.0059999988821064498,
    0.034714494024147946,
    -0.0089999988527065256,
    -0.069572895868103668,
    -0.030578415746029819,
    -0.04273869237059813,
    -0.0089999989133486835,
    -0.0089999989119302626,
    -0.041999998488987167,
    -0.041999998481231593,
    -0.0059999988940013793,
    -0.0059999988894703371,
    -0.0059999988705179419,
    -0.01199999878373248,
    -0.019999998678724928,
    -0.014999998732578845,
    -0.014377672160972743,
    -0.0059999988425225581
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
